# The deck's slide master ("theme1.xml") is themed with the "Integral"
# color scheme; the notes master ("theme2.xml") already carries the
# built-in Office color scheme. The commit swaps the two: the slide
# master (and therefore every slide) now takes on the standard Office
# color palette, while the notes master keeps the Integral palette.
#
# The font scheme (major/minor Latin typeface "Arial", and every other
# per-script fallback) and the fill/line/effect format scheme are
# byte-for-byte identical between the two themes, so the only
# observable difference to reproduce is the 12 theme colors.

function New-RGBVal($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$tcs = $p.SlideMaster.Theme.ThemeColorScheme

# Standard Office theme palette, in ThemeColorScheme index order:
# 1 Dark1, 2 Light1, 3 Dark2, 4 Light2,
# 5-10 Accent1..Accent6, 11 Hyperlink, 12 FollowedHyperlink
$tcs.Item(1).RGB  = New-RGBVal 0x00 0x00 0x00   # dk1      000000
$tcs.Item(2).RGB  = New-RGBVal 0xFF 0xFF 0xFF   # lt1      FFFFFF
$tcs.Item(3).RGB  = New-RGBVal 0x44 0x54 0x6A   # dk2      44546A
$tcs.Item(4).RGB  = New-RGBVal 0xE7 0xE6 0xE6   # lt2      E7E6E6
$tcs.Item(5).RGB  = New-RGBVal 0x5B 0x9B 0xD5   # accent1  5B9BD5
$tcs.Item(6).RGB  = New-RGBVal 0xED 0x7D 0x31   # accent2  ED7D31
$tcs.Item(7).RGB  = New-RGBVal 0xA5 0xA5 0xA5   # accent3  A5A5A5
$tcs.Item(8).RGB  = New-RGBVal 0xFF 0xC0 0x00   # accent4  FFC000
$tcs.Item(9).RGB  = New-RGBVal 0x44 0x72 0xC4   # accent5  4472C4
$tcs.Item(10).RGB = New-RGBVal 0x70 0xAD 0x47   # accent6  70AD47
$tcs.Item(11).RGB = New-RGBVal 0x05 0x63 0xC1   # hlink    0563C1
$tcs.Item(12).RGB = New-RGBVal 0x95 0x4F 0x72   # folHlink 954F72
